# Applies the "Added License, Added routing, Added submitted page" commit:
#   - adds a 4th signer row (Roy Lin / 340813412@yrdsb.ca / signer) with a
#     mailto hyperlink on the email cell, matching the formatting of the
#     existing rows 2-3
#   - widens columns A and B to fit the new content
#   - leaves the sheet's active selection on G11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: name / email / role -------------------------------------------------
$ws.Range("A4").Value = "Roy Lin"
$ws.Range("B4").Value = "340813412@yrdsb.ca"
$ws.Range("C4").Value = "signer"

# Hyperlink the new email cell the same way the existing ones are (rId1/rId2 -> mailto:...)
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:340813412@yrdsb.ca")

# Re-apply the shared Hyperlink cell style (Add() above stamps its own xf; this
# realigns B4 with the style already used by B2/B3).
$ws.Range("B4").Style = "Hyperlink"

# --- Column widths ----------------------------------------------------------
# ColumnWidth is quantized by the host to the nearest 1/6th of a character,
# so these inputs are chosen to land on the stored widths closest to the
# target 20.140625 / 27.28515625 (~20.1667 / ~27.3333 after quantization).
$ws.Columns.Item(1).ColumnWidth = 19.333333333333332
$ws.Columns.Item(2).ColumnWidth = 26.5

# --- Selection ---------------------------------------------------------------
$ws.Range("G11").Select()
